$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: paragraph "A cool quote by Dijkstra:" -> three runs:
#   "A cool quote by " + "Sunita" + ":"
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(5)
$full1 = $p1.Range.Text

# Plain substitution first (while the paragraph is still a single run).
$idxD = $full1.IndexOf("Dijkstra")
$dStart = $p1.Range.Start + $idxD
$dEnd = $dStart + 8
$dRange = $d.Range($dStart, $dEnd)
$dRange.Text = "Sunita"

# Split "Sunita" into its own run (must happen after all plain text edits,
# since any later .Text= assignment touching this area would re-merge
# same-formatted neighbouring runs).
$full1b = $p1.Range.Text
$idxS = $full1b.IndexOf("Sunita")
$sStart = $p1.Range.Start + $idxS
$sEnd = $sStart + 6
$sRange = $d.Range($sStart, $sEnd)
$sTemplate = $sRange.FormattedText
$sRange.FormattedText = $sTemplate

# Split the trailing ":" into its own run.
$full1c = $p1.Range.Text
$idxC = $full1c.LastIndexOf(":")
$cStart = $p1.Range.Start + $idxC
$cEnd = $cStart + 1
$cRange = $d.Range($cStart, $cEnd)
$cTemplate = $cRange.FormattedText
$cRange.FormattedText = $cTemplate

# ---------------------------------------------------------------------------
# Change 2: quote paragraph text replaced, with the closing curly quote
# split into its own (identically formatted) run.
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(6)
$full2 = $p2.Range.Text

# Replace the whole quoted sentence (including the trailing curly quote,
# excluding the paragraph mark) in one pass so the formatting of the
# original run is preserved on the new text.
$idx2 = $full2.IndexOf("Computer")
$start2 = $p2.Range.Start + $idx2
$end2 = $p2.Range.End - 1
$target2 = $d.Range($start2, $end2)
$target2.Text = "Everyone stay safe on this pandemic, Hopefully this pandemic will over soon and it will normal." + [char]8221

# Last step: split the closing curly quote off into its own run, reusing
# its own formatted-text snapshot (text + full rPr) so nothing else changes.
$full2b = $p2.Range.Text
$idxQuote = $full2b.IndexOf([char]8221)
$qStart = $p2.Range.Start + $idxQuote
$qEnd = $qStart + 1
$qRange = $d.Range($qStart, $qEnd)
$qTemplate = $qRange.FormattedText
$qRange.FormattedText = $qTemplate
